# Auto-generated Excel COM-interop script applying scheduled-runner price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Zeromus_Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 536.06665
$ws.Range("I33").Value2 = 567.29266
$ws.Range("K33").Value2 = 567.29266
$ws.Range("M33").Value2 = -338.29266

$ws.Range("H92").Value2 = 940.6923
$ws.Range("I92").Value2 = 747.1111
$ws.Range("J92").Value2 = 1376.25
$ws.Range("K92").Value2 = 747.1111
$ws.Range("L92").Value2 = 1376.25
$ws.Range("M92").Value2 = 500.8889
$ws.Range("N92").Value2 = -3872.25

$ws.Range("H111").Value2 = 2024.92
$ws.Range("I111").Value2 = 1962.15
$ws.Range("J111").Value2 = 2276
$ws.Range("K111").Value2 = 5886.450000000001
$ws.Range("L111").Value2 = 6828
$ws.Range("M111").Value2 = -2819.450000000001
$ws.Range("N111").Value2 = -12962

$ws.Range("H129").Value2 = 927.4
$ws.Range("J129").Value2 = 1040.88
$ws.Range("L129").Value2 = 3122.64
$ws.Range("N129").Value2 = -13122.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 835.7619
$ws.Range("I2").Value2 = 840.05884
$ws.Range("J2").Value2 = 817.5
$ws.Range("K2").Value2 = 840.05884
$ws.Range("L2").Value2 = 817.5
$ws.Range("M2").Value2 = -727.05884
$ws.Range("N2").Value2 = -1043.5

$ws.Range("H32").Value2 = 16365.032
$ws.Range("I32").Value2 = 3882
$ws.Range("J32").Value2 = 43913.793
$ws.Range("K32").Value2 = 3882
$ws.Range("L32").Value2 = 43913.793
$ws.Range("M32").Value2 = -3595
$ws.Range("N32").Value2 = -44487.793

$ws.Range("H112").Value2 = 24758
$ws.Range("I112").Value2 = 0
$ws.Range("J112").Value2 = 24758
$ws.Range("K112").Value2 = 0
$ws.Range("L112").ClearContents()
$ws.Range("M112").Value2 = 24758
$ws.Range("N112").Value2 = -27712

$ws.Range("H116").Value2 = 835.7619
$ws.Range("I116").Value2 = 840.05884
$ws.Range("J116").Value2 = 817.5
$ws.Range("K116").Value2 = 840.05884
$ws.Range("L116").Value2 = 817.5
$ws.Range("M116").Value2 = 1453.94116
$ws.Range("N116").Value2 = -5405.5

$ws.Range("H122").Value2 = 2536.1428
$ws.Range("I122").Value2 = 2656.25
$ws.Range("J122").Value2 = 2376
$ws.Range("K122").Value2 = 7968.75
$ws.Range("L122").Value2 = 7128
$ws.Range("M122").Value2 = -5518.75
$ws.Range("N122").Value2 = -12028

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 835.7619
$ws.Range("I3").Value2 = 840.05884
$ws.Range("J3").Value2 = 817.5
$ws.Range("K3").Value2 = 840.05884
$ws.Range("L3").Value2 = 817.5
$ws.Range("M3").Value2 = -726.05884
$ws.Range("N3").Value2 = -1045.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3228182.2
$ws.Range("I31").Value2 = 8697449
$ws.Range("K31").Value2 = 8697449
$ws.Range("M31").Value2 = -8697154

$ws.Range("H34").Value2 = 3228182.2
$ws.Range("I34").Value2 = 8697449
$ws.Range("K34").Value2 = 8697449
$ws.Range("M34").Value2 = -8697247

$ws.Range("H132").Value2 = 1880.1
$ws.Range("I132").Value2 = 1265.55
$ws.Range("J132").Value2 = 3109.2
$ws.Range("K132").Value2 = 3796.65
$ws.Range("L132").Value2 = 9327.599999999999
$ws.Range("M132").Value2 = -1266.65
$ws.Range("N132").Value2 = -14387.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value2 = 7114.2856
$ws.Range("I69").Value2 = 5500
$ws.Range("J69").Value2 = 7760
$ws.Range("K69").Value2 = 16500
$ws.Range("L69").Value2 = 23280
$ws.Range("M69").Value2 = -15689
$ws.Range("N69").Value2 = -24902

$ws.Range("H72").Value2 = 7114.2856
$ws.Range("I72").Value2 = 5500
$ws.Range("J72").Value2 = 7760
$ws.Range("K72").Value2 = 49500
$ws.Range("L72").Value2 = 69840
$ws.Range("M72").Value2 = -45444
$ws.Range("N72").Value2 = -77952

$ws.Range("H75").Value2 = 3557
$ws.Range("J75").Value2 = 4071.6667
$ws.Range("L75").Value2 = 12215.0001
$ws.Range("N75").Value2 = -14211.0001

$ws.Range("H78").Value2 = 3557
$ws.Range("J78").Value2 = 4071.6667
$ws.Range("L78").Value2 = 36645.0003
$ws.Range("N78").Value2 = -46629.0003

$ws.Range("H131").Value2 = 1852865.4
$ws.Range("I131").Value2 = 8333710.5
$ws.Range("J131").Value2 = 1195.3572
$ws.Range("K131").Value2 = 25001131.5
$ws.Range("L131").Value2 = 3586.0716
$ws.Range("M131").Value2 = -24996091.5
$ws.Range("N131").Value2 = -13666.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value2 = 5000
$ws.Range("J20").Value2 = 5000
$ws.Range("L20").Value2 = 5000
$ws.Range("N20").Value2 = -5490

$ws.Range("H25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("L25").ClearContents()
$ws.Range("N25").Value2 = 0

$ws.Range("H97").Value2 = 777.619
$ws.Range("I97").Value2 = 756.8946999999999
$ws.Range("K97").Value2 = 756.8946999999999
$ws.Range("M97").Value2 = -260.8946999999999

$ws.Range("H102").Value2 = 1339.5
$ws.Range("I102").Value2 = 1387
$ws.Range("J102").Value2 = 1007
$ws.Range("K102").Value2 = 1387
$ws.Range("L102").Value2 = 1007
$ws.Range("M102").Value2 = 235
$ws.Range("N102").Value2 = -4251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value2 = 5000
$ws.Range("I14").Value2 = 5000
$ws.Range("K14").Value2 = 5000
$ws.Range("M14").Value2 = -4828

$ws.Range("H16").Value2 = 934.6316
$ws.Range("I16").Value2 = 650.5333000000001
$ws.Range("J16").Value2 = 2000
$ws.Range("K16").Value2 = 650.5333000000001
$ws.Range("L16").Value2 = 2000
$ws.Range("M16").Value2 = -480.5333000000001
$ws.Range("N16").Value2 = -2340

$ws.Range("H46").Value2 = 1001.3
$ws.Range("I46").Value2 = 466.66666
$ws.Range("J46").Value2 = 1060.7037
$ws.Range("K46").Value2 = 466.66666
$ws.Range("L46").Value2 = 1060.7037
$ws.Range("M46").Value2 = -278.66666
$ws.Range("N46").Value2 = -1436.7037

$ws.Range("H55").Value2 = 265.10526
$ws.Range("I55").Value2 = 297.9091
$ws.Range("J55").Value2 = 220
$ws.Range("K55").Value2 = 297.9091
$ws.Range("L55").Value2 = 220
$ws.Range("M55").Value2 = -124.9091
$ws.Range("N55").Value2 = -566

$ws.Range("H110").Value2 = 31058
$ws.Range("J110").Value2 = 31058
$ws.Range("L110").Value2 = 31058
$ws.Range("N110").Value2 = -39238

$ws.Range("H132").Value2 = 7147777.5
$ws.Range("I132").Value2 = 13520831
$ws.Range("J132").Value2 = 2233.6667
$ws.Range("K132").Value2 = 40562493
$ws.Range("L132").Value2 = 6701.000100000001
$ws.Range("M132").Value2 = -40559963
$ws.Range("N132").Value2 = -11761.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value2 = 16606.25
$ws.Range("J86").Value2 = 17875
$ws.Range("L86").Value2 = 17875
$ws.Range("N86").Value2 = -20121

$ws.Range("H89").Value2 = 16606.25
$ws.Range("J89").Value2 = 17875
$ws.Range("L89").Value2 = 89375
$ws.Range("N89").Value2 = -100607

$ws.Range("H132").Value2 = 1211.8077
$ws.Range("I132").Value2 = 749.5806
$ws.Range("J132").Value2 = 1894.1428
$ws.Range("K132").Value2 = 2248.7418
$ws.Range("L132").Value2 = 5682.428400000001
$ws.Range("M132").Value2 = 281.2582000000002
$ws.Range("N132").Value2 = -10742.4284
